# Generate Report for Handback
# Marks the a.md / b.md localization files as handed-back for both the
# zh-cn and de-de target languages: updates the Overview status columns,
# fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on each language sheet (with a
# hyperlink + hyperlink styling on the new Target File cell, matching
# what the "a.md" source-file hyperlink already looks like), and widens
# a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b58606a87955b29669492dd45638ad7614be31d/e2e"
$handedBackStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both "zh-cn" and "de-de" status columns move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both
# tracked files (a.md row 2, b.md row 3). Those two columns also get
# wider to comfortably fit the longer status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $handedBackStatus
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------
# Per-language detail sheets ("zh-cn" and "de-de"): the handback just
# produced a target file (a.md, linked back to the source on GitHub)
# and a handback xliff, recorded with the moment it happened.
# ---------------------------------------------------------------------
function Update-LanguageSheet($SheetName, $HandbackFile, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the Status column (C) and the Latest Handback File column (J)
    # to fit their (now longer) contents.
    $ws.Columns.Item(3).ColumnWidth = 29.1666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.1666666666667

    # Status column: same "Ready for handoff" -> "Handed back: in sync
    # with en-US" wording update as the Overview sheet.
    $ws.Range("C2:C3").Value = $handedBackStatus

    # Latest Handback File / Latest Handback DateTime for both rows.
    $ws.Range("J2:J3").Value = $HandbackFile
    $ws.Range("K2:K3").Value = $HandbackDateTime

    # Latest Target File: the handed-back file content is a.md, linked
    # back to its GitHub source, same as column A's hyperlink.
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"

    # Rebuild hyperlinks in source order (A2, I2, A3, I3) so the
    # relationship ids land the same way Excel would assign them.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/a.md", $null, $null, "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), "$baseUrl/a.md", $null, $null, "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/b.md", $null, $null, "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), "$baseUrl/a.md", $null, $null, "a.md")
}

Update-LanguageSheet "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-26 14:47:49"
Update-LanguageSheet "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-26 14:47:57"
